# A new daily price record was reported for "Ajo" (Chino / Primera) at the
# Vega Monumental Concepción market, dated 2022-02-10 (Excel serial 44602).
# It belongs right after the existing row for 2021-10-15 (serial 44484),
# so insert a new row 18 and push every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value2 = 11
$ws.Range("B18").Value2 = "Vega Monumental Concepción"
$ws.Range("C18").Value2 = "Bíobío"
$ws.Range("D18").Value2 = 44602
$ws.Range("E18").Value2 = 8
$ws.Range("F18").Value2 = 100112003
$ws.Range("G18").Value2 = "Ajo"
$ws.Range("H18").Value2 = "Chino"
$ws.Range("I18").Value2 = "Primera"
$ws.Range("J18").Value2 = 170
$ws.Range("K18").Value2 = 20000
$ws.Range("L18").Value2 = 21000
$ws.Range("M18").Value2 = 20471
$ws.Range("N18").Value2 = '$/caja 10 kilos'
$ws.Range("O18").Value2 = "China"
$ws.Range("P18").Value2 = 2047
$ws.Range("Q18").Value2 = 10
$ws.Range("R18").Value2 = "Hortaliza"
